$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (name) and column G (mapbed) to fit the new ".17mb" suffix
$ws.Columns("B").ColumnWidth = 30.5
$ws.Columns("G").ColumnWidth = 38.33203125

# Row 2 (first sample) formulas get a ".17mb" tag inserted and are no longer
# shared with the rows below (edited independently from the fill below).
$ws.Range("B2").Formula = '=CONCATENATE(E2,".",C2,".17mb")'
$ws.Range("G2").Formula = '=CONCATENATE(E2,"_star.",C2,".17mb.bed.gz")'

# Rows 3-5 get the same ".17mb" tag, filled down from row 3 as a new shared
# formula group.
$ws.Range("B3").Formula = '=CONCATENATE(E3,".",C3,".17mb")'
$ws.Range("B4").Formula = '=CONCATENATE(E4,".",C4,".17mb")'
$ws.Range("B5").Formula = '=CONCATENATE(E5,".",C5,".17mb")'

$ws.Range("G3").Formula = '=CONCATENATE(E3,"_star.",C3,".17mb.bed.gz")'
$ws.Range("G4").Formula = '=CONCATENATE(E4,"_star.",C4,".17mb.bed.gz")'
$ws.Range("G5").Formula = '=CONCATENATE(E5,"_star.",C5,".17mb.bed.gz")'

# Selection moved to B3:B5 with active cell B3
$ws.Range("B3:B5").Select()

$ws.PageSetup.Orientation = 1
